# Apply the "Updated symbol list" crypto data refresh to Sheet1.
# For every changed cell: B/C (plain text: coin name / link) are written
# directly. D/E (price / volume%) look numeric, so Excel COM would silently
# re-type them as Number/Percent on a bare .Value assignment; forcing the
# range to Text first (NumberFormat "@") keeps them literal strings like the
# source file, then resetting .Style back to "Normal" drops the Text format
# again so no stray style survives the round-trip.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$r = $ws.Range("D2")
$r.NumberFormat = "@"
$r.Value = "329.55"
$r.Style = "Normal"
$r = $ws.Range("E2")
$r.NumberFormat = "@"
$r.Value = "6.58%"
$r.Style = "Normal"
# Row 3
$r = $ws.Range("D3")
$r.NumberFormat = "@"
$r.Value = "40.24"
$r.Style = "Normal"
$r = $ws.Range("E3")
$r.NumberFormat = "@"
$r.Value = "8.42%"
$r.Style = "Normal"
# Row 4
$r = $ws.Range("D4")
$r.NumberFormat = "@"
$r.Value = "5.598"
$r.Style = "Normal"
$r = $ws.Range("E4")
$r.NumberFormat = "@"
$r.Value = "9.12%"
$r.Style = "Normal"
# Row 5
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "0.08145"
$r.Style = "Normal"
$r = $ws.Range("E5")
$r.NumberFormat = "@"
$r.Value = "3.86%"
$r.Style = "Normal"
# Row 6
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = "4.551"
$r.Style = "Normal"
$r = $ws.Range("E6")
$r.NumberFormat = "@"
$r.Value = "3.33%"
$r.Style = "Normal"
# Row 7
$r = $ws.Range("D7")
$r.NumberFormat = "@"
$r.Value = "8.681"
$r.Style = "Normal"
$r = $ws.Range("E7")
$r.NumberFormat = "@"
$r.Value = "4.89%"
$r.Style = "Normal"
# Row 8
$r = $ws.Range("E8")
$r.NumberFormat = "@"
$r.Value = "5.29%"
$r.Style = "Normal"
# Row 10
$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = "0.9497"
$r.Style = "Normal"
$r = $ws.Range("E10")
$r.NumberFormat = "@"
$r.Value = "2.53%"
$r.Style = "Normal"
# Row 11
$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = "0.1274"
$r.Style = "Normal"
$r = $ws.Range("E11")
$r.NumberFormat = "@"
$r.Value = "8.88%"
$r.Style = "Normal"
# Row 12
$r = $ws.Range("D12")
$r.NumberFormat = "@"
$r.Value = "0.1989"
$r.Style = "Normal"
$r = $ws.Range("E12")
$r.NumberFormat = "@"
$r.Value = "4.83%"
$r.Style = "Normal"
# Row 13
$r = $ws.Range("D13")
$r.NumberFormat = "@"
$r.Value = "0.09223"
$r.Style = "Normal"
$r = $ws.Range("E13")
$r.NumberFormat = "@"
$r.Value = "4.07%"
$r.Style = "Normal"
# Row 14
$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = "0.03565"
$r.Style = "Normal"
$r = $ws.Range("E14")
$r.NumberFormat = "@"
$r.Value = "7.53%"
$r.Style = "Normal"
# Row 15
$r = $ws.Range("D15")
$r.NumberFormat = "@"
$r.Value = "0.09621"
$r.Style = "Normal"
$r = $ws.Range("E15")
$r.NumberFormat = "@"
$r.Value = "0.18%"
$r.Style = "Normal"
# Row 16
$r = $ws.Range("D16")
$r.NumberFormat = "@"
$r.Value = "0.001312"
$r.Style = "Normal"
$r = $ws.Range("E16")
$r.NumberFormat = "@"
$r.Value = "-4.86%"
$r.Style = "Normal"
# Row 17
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = "0.04433"
$r.Style = "Normal"
$r = $ws.Range("E17")
$r.NumberFormat = "@"
$r.Value = "2.09%"
$r.Style = "Normal"
# Row 18
$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$r = $ws.Range("D18")
$r.NumberFormat = "@"
$r.Value = "0.006083"
$r.Style = "Normal"
$r = $ws.Range("E18")
$r.NumberFormat = "@"
$r.Value = "-1.95%"
$r.Style = "Normal"
# Row 19
$ws.Range("B19").Value = "LEO"
$ws.Range("C19").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = "3.374"
$r.Style = "Normal"
$r = $ws.Range("E19")
$r.NumberFormat = "@"
$r.Value = "-0.43%"
$r.Style = "Normal"
# Row 20
$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = "0.3519"
$r.Style = "Normal"
$r = $ws.Range("E20")
$r.NumberFormat = "@"
$r.Value = "1.81%"
$r.Style = "Normal"
# Row 21
$ws.Range("B21").Value = "MCDex"
$ws.Range("C21").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = "7.476"
$r.Style = "Normal"
$r = $ws.Range("E21")
$r.NumberFormat = "@"
$r.Value = "17.26%"
$r.Style = "Normal"
# Row 22
$ws.Range("B22").Value = "ProBitToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = "0.1398"
$r.Style = "Normal"
$r = $ws.Range("E22")
$r.NumberFormat = "@"
$r.Value = "8.20%"
$r.Style = "Normal"
# Row 23
$ws.Range("B23").Value = "ZBToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = "0.2483"
$r.Style = "Normal"
$r = $ws.Range("E23")
$r.NumberFormat = "@"
$r.Value = "3.20%"
$r.Style = "Normal"
# Row 24
$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = "0.001258"
$r.Style = "Normal"
$r = $ws.Range("E24")
$r.NumberFormat = "@"
$r.Value = "4.64%"
$r.Style = "Normal"
# Row 25
$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = "0.004303"
$r.Style = "Normal"
$r = $ws.Range("E25")
$r.NumberFormat = "@"
$r.Value = "0.73%"
$r.Style = "Normal"
# Row 26
$r = $ws.Range("E26")
$r.NumberFormat = "@"
$r.Value = "-15.25%"
$r.Style = "Normal"
# Row 27
$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value = "0.0003985"
$r.Style = "Normal"
$r = $ws.Range("E27")
$r.NumberFormat = "@"
$r.Value = "37.21%"
$r.Style = "Normal"
# Row 39
$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = "0.02544"
$r.Style = "Normal"
$r = $ws.Range("E39")
$r.NumberFormat = "@"
$r.Value = "18.02%"
$r.Style = "Normal"
# Row 40
$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = "0.05224"
$r.Style = "Normal"
$r = $ws.Range("E40")
$r.NumberFormat = "@"
$r.Value = "4.33%"
$r.Style = "Normal"
# Row 41
$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = "0.007829"
$r.Style = "Normal"
$r = $ws.Range("E41")
$r.NumberFormat = "@"
$r.Value = "3.13%"
$r.Style = "Normal"
# Row 42
$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = "0.1436"
$r.Style = "Normal"
$r = $ws.Range("E42")
$r.NumberFormat = "@"
$r.Value = "5.99%"
$r.Style = "Normal"
# Row 43
$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = "0.008976"
$r.Style = "Normal"
$r = $ws.Range("E43")
$r.NumberFormat = "@"
$r.Value = "5.36%"
$r.Style = "Normal"
# Row 44
$r = $ws.Range("E44")
$r.NumberFormat = "@"
$r.Value = "8.69%"
$r.Style = "Normal"
# Row 45
$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = "0.009916"
$r.Style = "Normal"
$r = $ws.Range("E45")
$r.NumberFormat = "@"
$r.Value = "12.78%"
$r.Style = "Normal"
# Row 46
$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = "0.00006711"
$r.Style = "Normal"
$r = $ws.Range("E46")
$r.NumberFormat = "@"
$r.Value = "2.12%"
$r.Style = "Normal"
# Row 47
$r = $ws.Range("E47")
$r.NumberFormat = "@"
$r.Value = "-0.28%"
$r.Style = "Normal"
# Row 48
$r = $ws.Range("E48")
$r.NumberFormat = "@"
$r.Value = "-13.01%"
$r.Style = "Normal"
# Row 49
$r = $ws.Range("E49")
$r.NumberFormat = "@"
$r.Value = "59.14%"
$r.Style = "Normal"
# Row 50
$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = "0.00002097"
$r.Style = "Normal"
$r = $ws.Range("E50")
$r.NumberFormat = "@"
$r.Value = "-0.28%"
$r.Style = "Normal"
# Row 51
$r = $ws.Range("E51")
$r.NumberFormat = "@"
$r.Value = "-0.28%"
$r.Style = "Normal"
